# Repair the AlcoholUse mapping sheet:
# the "AlcoholUse.ObservationOfUse.StartDate" mapping value was sitting
# next to the wrong EHDS row (B13, "EHDSSubstanceUse.period"); it belongs
# next to "EHDSSubstanceUse.frequencyAndQuantity.period" (row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the misplaced value before wiping the cell.
$startDateValue = $ws.Range("B13").Value()

# Fully clear B13 (contents + formatting) so the cell disappears again.
$ws.Range("B13").Clear()

# Move the value to its correct row, matching the look of the sibling
# mapping cells in column B (Calibri, like B15/B19/...).
$ws.Range("B16").Value = $startDateValue
$ws.Range("B16").Font.Name = "Calibri"

# Leave the selection on the cell that was just fixed.
$ws.Range("B16").Select()
